$wb = $excel.ActiveWorkbook

# --- all_products: add form_factor column (P) ---
$apws = $wb.Worksheets.Item("all_products")
$apws.Range("P1").Value = "form_factor"
$apws.Range("P5").Value = "cardboard box"

# --- scif: add form_factor column (V), pulled via VLOOKUP from all_products ---
$scws = $wb.Worksheets.Item("scif")
$scws.Range("V1").Value = "form_factor"
for ($r = 2; $r -le 12; $r++) {
    $scws.Range("V$r").Formula = "=VLOOKUP(`$B$r, all_products!`$A`$2:`$P`$6, 16, 0)"
}

# --- matches: extend AutoFilter range to cover the new columns O:P ---
$mws = $wb.Worksheets.Item("matches")
$mws.AutoFilterMode = $false
$null = $mws.Range("A1:P52").AutoFilter()

# --- scif: extend AutoFilter range to cover the new column V ---
$scws.AutoFilterMode = $false
$null = $scws.Range("A1:V12").AutoFilter()

# --- selections (matches the recorded UI state after the edit) ---
$mws.Activate()
$mws.Range("I6").Select()

$apws.Activate()
$apws.Range("P2").Select()

$scws.Activate()
$scws.Range("V6").Select()

Write-Output "done"
